# Update "exp_temp_resp" (column H) values in sheet "s" to the latest
# recommended value, and make the "s" sheet the active/selected tab
# (previously the "info" sheet was active).

$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item("info")
$wsS    = $wb.Worksheets.Item("s")

# Update the exp_temp_resp values (column H, rows 2-4) on sheet "s"
$wsS.Range("H2").Value = 0.0948
$wsS.Range("H3").Value = 0.0948
$wsS.Range("H4").Value = 0.0948

# Reset selection on "info" sheet to B1 (was E18) while it's still active
$wsInfo.Activate()
$wsInfo.Range("B1").Select()

# Make "s" the active sheet/tab, with selection at A1 (unchanged from before)
$wsS.Activate()
$wsS.Range("A1").Select()
